$d = $word.ActiveDocument

# The document contains several "<id>...</id>" fragments that are split
# across three separate runs (Courier New "<id>", Arial "<name>", Courier
# New "</id>"). Re-running Find/Replace over each exact match with the
# same replacement text causes Word to re-type the match as a single run
# (inheriting the formatting of the first run of the match), which merges
# the three runs into one - matching the target edit.

$ids = @("p066r_5", "p066v_1", "p066v_2", "p066v_3", "p066v_4")

foreach ($id in $ids) {
    $needle = "<id>" + $id + "</id>"
    $d.Content.Find.Execute($needle, $true, $true, $false, $false, $false,
                             $true, 1, $false, $needle, 2) | Out-Null
}
